$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new data row (row 46) mirroring the existing rows' layout.
$row = 46

$ws.Cells.Item($row, 1).Value = 7
$ws.Cells.Item($row, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value = "Ñuble"

$ws.Cells.Item($row, 4).Value = 44448
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(45, 4).NumberFormat

$ws.Cells.Item($row, 5).Value = 16
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100108
$ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($row, 9).Value = 100108002
$ws.Cells.Item($row, 10).Value = "Mango"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 60
$ws.Cells.Item($row, 14).Value = 8500
$ws.Cells.Item($row, 15).Value = 9000
$ws.Cells.Item($row, 16).Value = 8750
$ws.Cells.Item($row, 17).Value = "$/bandeja 4 kilos"
$ws.Cells.Item($row, 18).Value = "Brasil"
$ws.Cells.Item($row, 19).Value = 2188
$ws.Cells.Item($row, 20).Value = 4
